$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Última Atualização" timestamp in column G for rows 2 through 50
for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 7).Value = "09/08/2025 09:03:26"
}
